# Updates Price (D) and Volume(1h) (E) columns for the cryptos list.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Reference style (the workbook default / unstyled cell) so that forcing
# column D to text ("@") for assignment does not leave a stray explicit
# style on the cell once the value has been written.
$refStyle = $ws.Range("B2").Style

# Map of cell -> new value. Column D values are plain numeric-looking
# strings (e.g. "1.002", "26.857.41") that must stay text, exactly as
# they were authored (matching the thousands-dot formatting already used
# in the sheet), so they are written through the text-coercion helper
# below. Column E values are already non-numeric (padded % strings) and
# are written directly.
$updates = [ordered]@{
    'D2' = '26.857.41'
    'E2' = '  -1.31%  '
    'D3' = '1.874.61'
    'D4' = '1.002'
    'E4' = '  -0.13%  '
    'D5' = '301.47'
    'E5' = '  -2.05%  '
    'E6' = '  -0.14%  '
    'D7' = '0.5359'
    'E7' = '  +1.85%  '
    'D8' = '0.3756'
    'E8' = '  -1.86%  '
    'D9' = '0.07186'
    'E9' = '  -1.58%  '
    'D10' = '21.57'
    'E10' = '  +0.11%  '
    'D11' = '0.8895'
    'E11' = '  -1.72%  '
    'D12' = '0.08145'
    'E12' = '  +0.72%  '
    'D13' = '1.887.70'
    'E13' = '  +2.78%  '
    'D14' = '93.48'
    'E14' = '  -2.47%  '
    'D15' = '5.319'
    'E15' = '  -0.89%  '
    'E16' = '  -0.08%  '
    'E17' = '  +0.92%  '
    'D18' = '0.000008544'
    'E18' = '  -1.54%  '
    'E19' = '  -0.12%  '
    'D20' = '26.891.34'
    'E20' = '  -1.32%  '
    'D21' = '4.986'
    'E21' = '  -2.64%  '
    'D22' = '10.63'
    'E22' = '  -1.88%  '
    'D23' = '6.402'
    'E23' = '  -1.44%  '
    'D24' = '2.300'
    'E24' = '  -1.66%  '
    'D25' = '146.17'
    'E25' = '  -2.57%  '
    'D26' = '18.07'
    'E26' = '  -0.99%  '
    'D27' = '1.731'
    'E27' = '  -0.64%  '
    'D28' = '114.07'
    'E28' = '  -2.30%  '
    'D29' = '4.721'
    'E29' = '  -2.53%  '
    'D30' = '4.615'
    'E30' = '  -5.32%  '
    'D31' = '0.09149'
    'E31' = '  -0.86%  '
    'D32' = '0.8139'
    'E32' = '  -0.26%  '
    'D33' = '0.05014'
    'E33' = '  -1.11%  '
    'D34' = '1.175'
    'E34' = '  -4.39%  '
    'D35' = '2.947'
    'E35' = '  -1.79%  '
    'D36' = '0.6040'
    'E36' = '  +5.45%  '
    'D37' = '3.215'
    'E37' = '  -4.41%  '
    'D38' = '2.611'
    'E38' = '  -3.66%  '
    'D39' = '0.01954'
    'E39' = '  -2.05%  '
    'E40' = '  -1.32%  '
    'D41' = '6.632'
    'E41' = '  +0.07%  '
    'D42' = '8.930'
    'E42' = '  -0.89%  '
    'D43' = '115.13'
    'E43' = '  -1.58%  '
    'D44' = '0.5087'
    'E44' = '  +3.25%  '
    'D45' = '0.1493'
    'E45' = '  -1.92%  '
    'E46' = '  -0.15%  '
    'D47' = '9.943'
    'E47' = '  -2.03%  '
    'E48' = '  -0.21%  '
    'D49' = '37.70'
    'E49' = '  -2.19%  '
    'D50' = '0.06052'
    'E50' = '  +1.43%  '
    'D51' = '62.25'
    'E51' = '  -3.02%  '
}

foreach ($addr in $updates.Keys) {
    $value = $updates[$addr]
    $cell = $ws.Range($addr)
    if ($addr[0] -eq "D") {
        # Force text so numeric-looking strings (e.g. "1.002") are not
        # reinterpreted as numbers, then restore the original (default)
        # style so no formatting diff is introduced.
        $cell.NumberFormat = "@"
        $cell.Value = $value
        $cell.Style = $refStyle
    } else {
        $cell.Value = $value
    }
}
